$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Status for the first user story row to "In Progress"
$ws.Range("D2").Value = "In Progress"

# Update the active cell selection to F1 (as left by the author after editing)
$ws.Range("F1").Select()
